$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New radius (nm) values for B3:B10 (B2 and B11 stay the same)
$ws.Range("B3").Value = 61.11111111111111
$ws.Range("B4").Value = 72.22222222222223
$ws.Range("B5").Value = 83.33333333333333
$ws.Range("B6").Value = 94.44444444444444
$ws.Range("B7").Value = 105.5555555555556
$ws.Range("B8").Value = 116.6666666666667
$ws.Range("B9").Value = 127.7777777777778
$ws.Range("B10").Value = 138.8888888888889

# New pixel size (nm) values for E2:E11 (all rows 2-11 change from 2 to 2.5)
$ws.Range("E2:E11").Value = 2.5
